# Insert a new data row at row 164 (pushing existing rows 164-192 down to 165-193)
# to add the latest week's price record for Ciboulette at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 164; this shifts formatting (including the
# date number format on column D) down along with the existing data.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record's values.
$ws.Cells.Item(164, 1).Value = 4
$ws.Cells.Item(164, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(164, 3).Value = "Los Lagos"
$ws.Cells.Item(164, 4).Value = 44637
$ws.Cells.Item(164, 5).Value = 10
$ws.Cells.Item(164, 6).Value = 100112039
$ws.Cells.Item(164, 7).Value = "Ciboulette"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 80
$ws.Cells.Item(164, 11).Value = 4000
$ws.Cells.Item(164, 12).Value = 4000
$ws.Cells.Item(164, 13).Value = 4000
$ws.Cells.Item(164, 14).Value = "$/docena de atados"
$ws.Cells.Item(164, 15).Value = "Región Metropolitana"
$ws.Cells.Item(164, 16).Value = 1333
$ws.Cells.Item(164, 17).Value = 3
$ws.Cells.Item(164, 18).Value = "Hortaliza"
